$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Salary" header column right after "Year Exp" (column O),
# matching the formatting already used by the rest of the header row.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = "Salary"
$excel.CutCopyMode = $false
